# Delete row 21 ("「いのち」..." entry) so all following rows shift up by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(21).Delete()
